# Final(?) push of Excel templates for testing.
# The "prodprice" tag row (B33) is no longer needed, so its value is cleared.
# Since "prodprice" was the only usage of that shared string, Excel drops it
# from the shared-string table on save, shifting the indices of every other
# string that sorted after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the cell that held "prodprice" (row 33, column B).
$ws.Range("B33").ClearContents()

# Reflect where the user ended up after the edit: selection on B33.
$ws.Range("B33").Select()
